$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B-E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values for columns B-E
$ws.Range("B2").Value = 5.4999095751223175
$ws.Range("C2").Value = 1.2637975731455251
$ws.Range("D2").Value = 1.783092086562087
$ws.Range("E2").Value = 1.7468665119103974

# Update row 3 values for columns B-E
$ws.Range("B3").Value = 4.2882762550519846
$ws.Range("C3").Value = 7.3636155101970564
$ws.Range("D3").Value = 9.4515205349522233
$ws.Range("E3").Value = -0.81485189840192296

# Update the selection to match the new range B1:E3
$ws.Range("B1:E3").Select()
